$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve formatting of column D while writing text-like numeric values
$colD = $ws.Range("D2:D51")
$origStyle = $colD.Style
$colD.NumberFormat = "@"

$ws.Range("D2").Value = "25.684.87"
$ws.Range("E2").Value = "  -3.49%  "
$ws.Range("D3").Value = "1.742.24"
$ws.Range("E3").Value = "  -5.80%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "235.94"
$ws.Range("E5").Value = "  -10.22%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "0.4921"
$ws.Range("E7").Value = "  -7.97%  "
$ws.Range("D8").Value = "41.43"
$ws.Range("E8").Value = "  -8.36%  "
$ws.Range("D9").Value = "0.2544"
$ws.Range("E9").Value = "  -19.18%  "
$ws.Range("D10").Value = "0.06018"
$ws.Range("E10").Value = "  -13.04%  "
$ws.Range("D11").Value = "1.743.56"
$ws.Range("E11").Value = "  -5.60%  "
$ws.Range("D12").Value = "0.06826"
$ws.Range("D13").Value = "14.82"
$ws.Range("E13").Value = "  -21.27%  "
$ws.Range("D14").Value = "4.439"
$ws.Range("E14").Value = "  -12.11%  "
$ws.Range("D15").Value = "76.40"
$ws.Range("E15").Value = "  -14.76%  "
$ws.Range("D16").Value = "0.5685"
$ws.Range("E16").Value = "  -26.10%  "
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "25.717.80"
$ws.Range("E19").Value = "  -3.45%  "
$ws.Range("E20").Value = "  -19.89%  "
$ws.Range("D21").Value = "0.000006548"
$ws.Range("E21").Value = "  -17.75%  "
$ws.Range("D22").Value = "1.967.16"
$ws.Range("E22").Value = "  -5.81%  "
$ws.Range("D23").Value = "4.016"
$ws.Range("E23").Value = "  -13.56%  "
$ws.Range("D24").Value = "5.050"
$ws.Range("E24").Value = "  -16.10%  "
$ws.Range("D25").Value = "7.933"
$ws.Range("E25").Value = "  -15.10%  "
$ws.Range("D26").Value = "136.98"
$ws.Range("E26").Value = "  -3.13%  "
$ws.Range("D27").Value = "1.473"
$ws.Range("E27").Value = "  -12.89%  "
$ws.Range("D28").Value = "1.815"
$ws.Range("E28").Value = "  -18.11%  "
$ws.Range("D29").Value = "14.65"
$ws.Range("E29").Value = "  -13.99%  "
$ws.Range("D30").Value = "101.73"
$ws.Range("E30").Value = "  -8.77%  "
$ws.Range("D31").Value = "3.744"
$ws.Range("E31").Value = "  -13.00%  "
$ws.Range("D32").Value = "0.07964"
$ws.Range("E32").Value = "  -9.29%  "
$ws.Range("D33").Value = "3.393"
$ws.Range("E33").Value = "  -17.48%  "
$ws.Range("D34").Value = "0.04382"
$ws.Range("E34").Value = "  -9.73%  "
$ws.Range("D35").Value = "1.000"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "2.613"
$ws.Range("E36").Value = "  -10.28%  "
$ws.Range("D37").Value = "0.9810"
$ws.Range("E37").Value = "  -13.86%  "
$ws.Range("D38").Value = "0.5966"
$ws.Range("E38").Value = "  -19.20%  "
$ws.Range("D39").Value = "2.659"
$ws.Range("E39").Value = "  -14.62%  "
$ws.Range("D40").Value = "1.911"
$ws.Range("E40").Value = "  -18.08%  "
$ws.Range("D41").Value = "1.002"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "0.01513"
$ws.Range("E42").Value = "  -12.69%  "
$ws.Range("D43").Value = "101.59"
$ws.Range("E43").Value = "  -6.39%  "
$ws.Range("D44").Value = "0.7529"
$ws.Range("E44").Value = "  -16.85%  "
$ws.Range("D45").Value = "5.157"
$ws.Range("E45").Value = "  -12.64%  "
$ws.Range("D46").Value = "0.3744"
$ws.Range("E46").Value = "  -22.47%  "
$ws.Range("D47").Value = "0.05227"
$ws.Range("E47").Value = "  -10.11%  "
$ws.Range("D48").Value = "0.1065"
$ws.Range("E48").Value = "  -14.67%  "
$ws.Range("D49").Value = "30.05"
$ws.Range("E49").Value = "  -14.28%  "
$ws.Range("D50").Value = "52.06"
$ws.Range("E50").Value = "  -13.77%  "
$ws.Range("D51").Value = "5.791"
$ws.Range("E51").Value = "  -24.66%  "

# Restore original column D formatting/style
$colD.Style = $origStyle
